$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.088.85"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.788.56"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'226.72"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'0.547"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'32.11"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "'0.0693"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "2.046.45"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'11.55"
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("D14").Value = "1.790.49"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "'0.622"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "34.079.33"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'4.20"
$ws.Range("D18").Value = "'68.00"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'244.61"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "'10.92"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'4.10"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "'162.59"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'7.19"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "'16.29"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "'0.0521"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "1.416.05"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  +5.69%  "
$ws.Range("D40").Value = "'80.52"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "'13.39"
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "'1.07"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -5.78%  "
$ws.Range("D49").Value = "'107.19"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "1.947.76"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.23%  "
